$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2 = @{ "B"="25.93855173630075"; "C"="7.281048931762184"; "D"="12.81472341064355"; "E"="12.46156716183474"; "G"="66.66898528846623"; "H"="23.6134587998259"; "J"="7.805370772747803"; "L"="12.90968755769887"; "M"="21.68533418851077"; "N"="22.56163371710123" }
  3 = @{ "B"="25.69182149556792"; "C"="6.909330383966485"; "D"="12.83149741986555"; "E"="12.48250287686244"; "G"="66.53820328437678"; "H"="23.64253773236256"; "J"="7.791883255362388"; "L"="12.92510626394327"; "M"="21.65132401001095"; "N"="22.63129885194341" }
  4 = @{ "B"="25.54586903253012"; "C"="6.67242734240467"; "D"="12.84377687596253"; "E"="12.49610575424739"; "G"="66.4751484562546"; "H"="23.6648869339852"; "J"="7.783383981621722"; "L"="12.93639489805079"; "M"="21.63462894771272"; "N"="22.67609245782819" }
  5 = @{ "B"="25.4878466240949"; "C"="6.573841302988503"; "D"="12.8492785794027"; "E"="12.50183770669634"; "G"="66.45379611712264"; "H"="23.67512220337397"; "J"="7.779865596793141"; "L"="12.94145328109762"; "M"="21.62888311933336"; "N"="22.69485507474165" }
  6 = @{ "B"="25.47830161649043"; "C"="6.557351908486068"; "D"="12.85022218714429"; "E"="12.5028009033015"; "G"="66.45051298784578"; "H"="23.67688980145436"; "J"="7.779278055962719"; "L"="12.94232089728915"; "M"="21.6279930216823"; "N"="22.69800136168992" }
  7 = @{ "B"="25.54508055478894"; "C"="6.671105869236051"; "D"="12.84384905890584"; "E"="12.49618229276344"; "G"="66.47484290141603"; "H"="23.6650204075479"; "J"="7.783336753681207"; "L"="12.93646126191271"; "M"="21.63454717003315"; "N"="22.67634343517949" }
  8 = @{ "B"="25.85236667708597"; "C"="7.154752354197578"; "D"="12.82009601058805"; "E"="12.46863080283379"; "G"="66.62031738203943"; "H"="23.62255097904882"; "J"="7.800765452695443"; "L"="12.91462598171944"; "M"="21.67274138507576"; "N"="22.58523594082799" }
  9 = @{ "B"="26.4961007965651"; "C"="8.029599759403579"; "D"="12.78923505594732"; "E"="12.42051630621096"; "G"="67.04196499098084"; "H"="23.57503043448321"; "J"="7.833218481896048"; "L"="12.88625216441348"; "M"="21.78063160020079"; "N"="22.42253857663123" }
  10 = @{ "B"="26.99020832511087"; "C"="8.62227382650571"; "D"="12.77615115256752"; "E"="12.38873951543222"; "G"="67.4339962553161"; "H"="23.56204037163194"; "J"="7.856017517743505"; "L"="12.87420004191777"; "M"="21.87965449413405"; "N"="22.31265810248163" }
  11 = @{ "B"="27.21871711683113"; "C"="8.880216812965665"; "D"="12.77228114155153"; "E"="12.37505241127177"; "G"="67.62992172715278"; "H"="23.56091124583913"; "J"="7.866164206158206"; "L"="12.87062292472435"; "M"="21.9289009025211"; "N"="22.26475006883747" }
  12 = @{ "B"="27.30571090109926"; "C"="8.976163782881272"; "D"="12.77111484806308"; "E"="12.36997942560767"; "G"="67.70661223283494"; "H"="23.56117205611853"; "J"="7.869974261562307"; "L"="12.86954188641384"; "M"="21.94814424557268"; "N"="22.24690605083865" }
  13 = @{ "B"="27.28695585804311"; "C"="8.955577583401377"; "D"="12.77135272734982"; "E"="12.37106709777552"; "G"="67.68998501449735"; "H"="23.5610852573758"; "J"="7.86915513284667"; "L"="12.86976254987262"; "M"="21.94397354791393"; "N"="22.2507358537515" }
  14 = @{ "B"="27.22586530860271"; "C"="8.888145375516649"; "D"="12.77217919494073"; "E"="12.37463285126596"; "G"="67.63618123648189"; "H"="23.56091890202758"; "J"="7.866478308349271"; "L"="12.87052850786803"; "M"="21.93047219864984"; "N"="22.26327606961401" }
  15 = @{ "B"="27.18850358523187"; "C"="8.846614451126493"; "D"="12.77272438727238"; "E"="12.37683129300534"; "G"="67.60354915125137"; "H"="23.56090667680081"; "J"="7.864834473922934"; "L"="12.87103328693824"; "M"="21.92227940675879"; "N"="22.2709960565268" }
  16 = @{ "B"="26.97534353502755"; "C"="8.605177135073832"; "D"="12.77644594709986"; "E"="12.38964942100547"; "G"="67.42154339841736"; "H"="23.56221044509468"; "J"="7.855349889360235"; "L"="12.87447211968044"; "M"="21.87651987702309"; "N"="22.31583071867147" }
  17 = @{ "B"="26.84548062804306"; "C"="8.454034135362457"; "D"="12.779262134432"; "E"="12.39770938345923"; "G"="67.31437255536025"; "H"="23.56423539561496"; "J"="7.849473878425939"; "L"="12.87706945929462"; "M"="21.84951738394438"; "N"="22.3438666369237" }
  18 = @{ "B"="26.77114315305297"; "C"="8.366004939066066"; "D"="12.78107790825644"; "E"="12.40241760554647"; "G"="67.25438838025821"; "H"="23.56584998237829"; "J"="7.84607299976596"; "L"="12.87874276333743"; "M"="21.83438222364431"; "N"="22.36018772343297" }
  19 = @{ "B"="26.74603722306474"; "C"="8.336013342132045"; "D"="12.78172636134343"; "E"="12.40402416801615"; "G"="67.23436436326998"; "H"="23.56647388605519"; "J"="7.844917886533058"; "L"="12.8793401373756"; "M"="21.82932599372523"; "N"="22.36574738191007" }
  20 = @{ "B"="26.85926842936557"; "C"="8.470237383749915"; "D"="12.7789420646054"; "E"="12.39684390262955"; "G"="67.32560972813934"; "H"="23.56397326619324"; "J"="7.850101578741691"; "L"="12.87677440574747"; "M"="21.85235092969391"; "N"="22.34086193230156" }
  21 = @{ "B"="27.24379712992094"; "C"="8.907999184566382"; "D"="12.77192832321869"; "E"="12.37358252065988"; "G"="67.65191719192192"; "H"="23.56094907568755"; "J"="7.86726543224855"; "L"="12.87029610776505"; "M"="21.93442180657604"; "N"="22.2595846320189" }
  22 = @{ "B"="27.49776975750614"; "C"="9.183997655818057"; "D"="12.76908822668843"; "E"="12.35902095484134"; "G"="67.87971953622437"; "H"="23.56298525515818"; "J"="7.878294893204271"; "L"="12.86765628622645"; "M"="21.9915220267736"; "N"="22.20820008847101" }
  23 = @{ "B"="27.36200105531922"; "C"="9.037631551462619"; "D"="12.77044456788231"; "E"="12.36673422720715"; "G"="67.75681801459072"; "H"="23.56153110444071"; "J"="7.872425462674261"; "L"="12.86891951975222"; "M"="21.96073300498399"; "N"="22.23546656263386" }
  24 = @{ "B"="26.85303395430309"; "C"="8.462915426608951"; "D"="12.77908615544459"; "E"="12.39723495460313"; "G"="67.32052432337476"; "H"="23.56409037193325"; "J"="7.849817866086685"; "L"="12.87690723856179"; "M"="21.85106867147666"; "N"="22.34221972724708" }
  25 = @{ "B"="26.31794719926099"; "C"="7.801407630033966"; "D"="12.79589988773857"; "E"="12.43290278925087"; "G"="66.91337310750069"; "H"="23.58404362717862"; "J"="7.8246229555873"; "L"="12.89238254295334"; "M"="21.74794622920801"; "N"="22.46485168076452" }
}

foreach ($r in $data.Keys) {
  $rowVals = $data[$r]
  foreach ($c in $rowVals.Keys) {
    $ws.Range("$c$r").Value = [double]$rowVals[$c]
  }
}

Write-Output "Updated $($data.Keys.Count) rows"